$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.583.70"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.842.02"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5273"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3165"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06798"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7835"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07797"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "1.835.75"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.013"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007932"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "26.614.73"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "2.074.88"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.613"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.335"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.221"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.682"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.212"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08701"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.080"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04861"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7322"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.139"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.865"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.095"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.345"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.10%  "
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4825"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9042"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.910"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.708"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4200"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.094"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1245"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8933"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
